$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Robo1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.541588333333333
$ws.Cells.Item(2, 8).Value = 7.624765
$ws.Cells.Item(2, 9).Value = 0.0417373532195736
$ws.Cells.Item(2, 10).Value = 0.0417373532195736
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.6629413333333334
$ws.Cells.Item(2, 14).Value = 1.988824
$ws.Cells.Item(2, 15).Value = 0.01938545156182102
$ws.Cells.Item(2, 16).Value = 0.01938545156182102
$ws.Cells.Item(2, 17).Value = 1.684923958484445
$ws.Cells.Item(2, 18).Value = 15.16431562636
$ws.Cells.Item(2, 19).Value = 0.0008090974391566586
$ws.Cells.Item(2, 20).Value = 0.0008090974391566586

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Robo1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.541588333333333
$ws.Cells.Item(3, 8).Value = 7.624765
$ws.Cells.Item(3, 9).Value = 0.0417373532195736
$ws.Cells.Item(3, 10).Value = 0.0417373532195736
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 29.811843
$ws.Cells.Item(3, 14).Value = 89.435529
$ws.Cells.Item(3, 15).Value = 0.8717453707996982
$ws.Cells.Item(3, 16).Value = 0.8717453707996982
$ws.Cells.Item(3, 17).Value = 75.769432363965
$ws.Cells.Item(3, 18).Value = 681.924891275685
$ws.Cells.Item(3, 19).Value = 0.03638434445859516
$ws.Cells.Item(3, 20).Value = 0.03638434445859516

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Robo1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.541588333333333
$ws.Cells.Item(4, 8).Value = 7.624765
$ws.Cells.Item(4, 9).Value = 0.0417373532195736
$ws.Cells.Item(4, 10).Value = 0.0417373532195736
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.723095
$ws.Cells.Item(4, 14).Value = 11.169285
$ws.Cells.Item(4, 15).Value = 0.1088691776384809
$ws.Cells.Item(4, 16).Value = 0.1088691776384809
$ws.Cells.Item(4, 17).Value = 9.462574815891667
$ws.Cells.Item(4, 18).Value = 85.16317334302501
$ws.Cells.Item(4, 19).Value = 0.00454391132182178
$ws.Cells.Item(4, 20).Value = 0.00454391132182178

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Robo1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 22.91769833333333
$ws.Cells.Item(5, 8).Value = 68.753095
$ws.Cells.Item(5, 9).Value = 0.376348938092374
$ws.Cells.Item(5, 10).Value = 0.376348938092374
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6629413333333334
$ws.Cells.Item(5, 14).Value = 1.988824
$ws.Cells.Item(5, 15).Value = 0.01938545156182102
$ws.Cells.Item(5, 16).Value = 0.01938545156182102
$ws.Cells.Item(5, 17).Value = 15.19308949003111
$ws.Cells.Item(5, 18).Value = 136.73780541028
$ws.Cells.Item(5, 19).Value = 0.007295694109732492
$ws.Cells.Item(5, 20).Value = 0.007295694109732492

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Robo1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 22.91769833333333
$ws.Cells.Item(6, 8).Value = 68.753095
$ws.Cells.Item(6, 9).Value = 0.376348938092374
$ws.Cells.Item(6, 10).Value = 0.376348938092374
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 29.811843
$ws.Cells.Item(6, 14).Value = 89.435529
$ws.Cells.Item(6, 15).Value = 0.8717453707996982
$ws.Cells.Item(6, 16).Value = 0.8717453707996982
$ws.Cells.Item(6, 17).Value = 683.218824634695
$ws.Cells.Item(6, 18).Value = 6148.969421712255
$ws.Cells.Item(6, 19).Value = 0.3280804445874092
$ws.Cells.Item(6, 20).Value = 0.3280804445874092

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Robo1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 22.91769833333333
$ws.Cells.Item(7, 8).Value = 68.753095
$ws.Cells.Item(7, 9).Value = 0.376348938092374
$ws.Cells.Item(7, 10).Value = 0.376348938092374
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.723095
$ws.Cells.Item(7, 14).Value = 11.169285
$ws.Cells.Item(7, 15).Value = 0.1088691776384809
$ws.Cells.Item(7, 16).Value = 0.1088691776384809
$ws.Cells.Item(7, 17).Value = 85.32476807634167
$ws.Cells.Item(7, 18).Value = 767.922912687075
$ws.Cells.Item(7, 19).Value = 0.0409727993952323
$ws.Cells.Item(7, 20).Value = 0.0409727993952323

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Ncam1"
$ws.Cells.Item(8, 3).Value = "Robo1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.04007766666666667
$ws.Cells.Item(8, 8).Value = 0.120233
$ws.Cells.Item(8, 9).Value = 0.0006581458169070119
$ws.Cells.Item(8, 10).Value = 0.0006581458169070119
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.6629413333333334
$ws.Cells.Item(8, 14).Value = 1.988824
$ws.Cells.Item(8, 15).Value = 0.01938545156182102
$ws.Cells.Item(8, 16).Value = 0.01938545156182102
$ws.Cells.Item(8, 17).Value = 0.02656914177688889
$ws.Cells.Item(8, 18).Value = 0.239122275992
$ws.Cells.Item(8, 19).Value = 0.00001275845385426601
$ws.Cells.Item(8, 20).Value = 0.00001275845385426601

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Ncam1"
$ws.Cells.Item(9, 3).Value = "Robo1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.04007766666666667
$ws.Cells.Item(9, 8).Value = 0.120233
$ws.Cells.Item(9, 9).Value = 0.0006581458169070119
$ws.Cells.Item(9, 10).Value = 0.0006581458169070119
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 29.811843
$ws.Cells.Item(9, 14).Value = 89.435529
$ws.Cells.Item(9, 15).Value = 0.8717453707996982
$ws.Cells.Item(9, 16).Value = 0.8717453707996982
$ws.Cells.Item(9, 17).Value = 1.194789106473
$ws.Cells.Item(9, 18).Value = 10.753101958257
$ws.Cells.Item(9, 19).Value = 0.0005737355691998733
$ws.Cells.Item(9, 20).Value = 0.0005737355691998733

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Ncam1"
$ws.Cells.Item(10, 3).Value = "Robo1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.04007766666666667
$ws.Cells.Item(10, 8).Value = 0.120233
$ws.Cells.Item(10, 9).Value = 0.0006581458169070119
$ws.Cells.Item(10, 10).Value = 0.0006581458169070119
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.723095
$ws.Cells.Item(10, 14).Value = 11.169285
$ws.Cells.Item(10, 15).Value = 0.1088691776384809
$ws.Cells.Item(10, 16).Value = 0.1088691776384809
$ws.Cells.Item(10, 17).Value = 0.1492129603783334
$ws.Cells.Item(10, 18).Value = 1.342916643405
$ws.Cells.Item(10, 19).Value = 0.00007165179385287258
$ws.Cells.Item(10, 20).Value = 0.00007165179385287258

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Ncam1"
$ws.Cells.Item(11, 3).Value = "Robo1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 35.395449
$ws.Cells.Item(11, 8).Value = 106.186347
$ws.Cells.Item(11, 9).Value = 0.5812555628711454
$ws.Cells.Item(11, 10).Value = 0.5812555628711454
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.6629413333333334
$ws.Cells.Item(11, 14).Value = 1.988824
$ws.Cells.Item(11, 15).Value = 0.01938545156182102
$ws.Cells.Item(11, 16).Value = 0.01938545156182102
$ws.Cells.Item(11, 17).Value = 23.465106153992
$ws.Cells.Item(11, 18).Value = 211.185955385928
$ws.Cells.Item(11, 19).Value = 0.0112679015590776
$ws.Cells.Item(11, 20).Value = 0.0112679015590776

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Ncam1"
$ws.Cells.Item(12, 3).Value = "Robo1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 35.395449
$ws.Cells.Item(12, 8).Value = 106.186347
$ws.Cells.Item(12, 9).Value = 0.5812555628711454
$ws.Cells.Item(12, 10).Value = 0.5812555628711454
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 29.811843
$ws.Cells.Item(12, 14).Value = 89.435529
$ws.Cells.Item(12, 15).Value = 0.8717453707996982
$ws.Cells.Item(12, 16).Value = 0.8717453707996982
$ws.Cells.Item(12, 17).Value = 1055.203568502507
$ws.Cells.Item(12, 18).Value = 9496.832116522563
$ws.Cells.Item(12, 19).Value = 0.5067068461844939
$ws.Cells.Item(12, 20).Value = 0.5067068461844939

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Ncam1"
$ws.Cells.Item(13, 3).Value = "Robo1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 35.395449
$ws.Cells.Item(13, 8).Value = 106.186347
$ws.Cells.Item(13, 9).Value = 0.5812555628711454
$ws.Cells.Item(13, 10).Value = 0.5812555628711454
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.723095
$ws.Cells.Item(13, 14).Value = 11.169285
$ws.Cells.Item(13, 15).Value = 0.1088691776384809
$ws.Cells.Item(13, 16).Value = 0.1088691776384809
$ws.Cells.Item(13, 17).Value = 131.780619194655
$ws.Cells.Item(13, 18).Value = 1186.025572751895
$ws.Cells.Item(13, 19).Value = 0.06328081512757391
$ws.Cells.Item(13, 20).Value = 0.06328081512757391
